# Apply the commit's changes to the "按三大门类分工业产能利用率" workbook:
#   1. Delete columns R:U entirely (the duplicate/derived
#      工业产能利用率_制造业 / 电力... / 工业产能利用率 / 采矿业 columns),
#      shifting the dimension from A1:U5 down to A1:Q5.
#   2. Swap the data rows for "2021年B" and "2021年C" (rows 3 and 4),
#      so the "2021年C" figures now sit in row 3 and the "2021年B"
#      figures sit in row 4 (columns A:Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove columns R:U, shifting everything left ---------------------
$ws.Range("R1:U5").Delete(-4161) | Out-Null   # -4161 = xlShiftToLeft

# --- 2. Swap row 3 ("2021年B") and row 4 ("2021年C") across A:Q ----------
$lastCol = 17   # column Q

for ($c = 1; $c -le $lastCol; $c++) {
    $cell3 = $ws.Cells.Item(3, $c)
    $cell4 = $ws.Cells.Item(4, $c)

    $v3 = $cell3.Value2
    $v4 = $cell4.Value2

    $cell3.Value = $v4
    $cell4.Value = $v3
}
